$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E8: was numeric 5 -> becomes text "4 a 5.5" (the "Alcool (%)" range label)
$ws.Range("E8").Value = "4 a 5.5"

# E14: was numeric 13 -> becomes text "13 a 15" (new range label)
$ws.Range("E14").Value = "13 a 15"

# Move the active selection from E17 to E9
$ws.Range("E9").Select()
